$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.135.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.898.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.94%  "

$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4618"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3901"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07894"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9908"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.880.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.067"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.54%  "

$ws.Range("E14").Value = "  -0.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06992"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.52%  "

$ws.Range("E17").Value = "  -0.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009990"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.148.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.320"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.105.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.111"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.960"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "118.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.881"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09333"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.255"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.326"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.149"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05794"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.28%  "

$ws.Range("E37").Value = "  -0.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02085"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.001"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.733"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5694"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.89%  "

$ws.Range("E42").Value = "  -0.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.717"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5355"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.172"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07015"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.855"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.558"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.044"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.19%  "
